# Update the drug "stock" values (column D) for amoxicillin, chloramphenicol,
# polymyxinB and tetracycline.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 32
$ws.Range("D3").Value = 32
$ws.Range("D6").Value = 16
$ws.Range("D7").Value = 16

# Widen column E (micRef) a bit.
$ws.Columns.Item(5).ColumnWidth = 18.33

# Move the active selection from F8 to D8.
$ws.Range("D8").Select()
